$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on price cells whose new values would otherwise
# be auto-converted to numbers (losing significant trailing zeros / exact text).
$ws.Range("D5:D7").NumberFormat = "@"
$ws.Range("D9:D12").NumberFormat = "@"
$ws.Range("D15:D16").NumberFormat = "@"
$ws.Range("D20:D24").NumberFormat = "@"
$ws.Range("D26:D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D36:D38").NumberFormat = "@"
$ws.Range("D40:D41").NumberFormat = "@"
$ws.Range("D46:D48").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated crypto price/volume values from the refreshed feed.
$ws.Range("D2").Value = "44.225.06"
$ws.Range("E2").Value = "  +1.87%  "
$ws.Range("D3").Value = "2.362.25"
$ws.Range("E3").Value = "  -0.51%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "244.34"
$ws.Range("E5").Value = "  +3.58%  "
$ws.Range("D6").Value = "0.681"
$ws.Range("E6").Value = "  +5.27%  "
$ws.Range("D7").Value = "74.56"
$ws.Range("E7").Value = "  +4.69%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "0.580"
$ws.Range("E9").Value = "  +25.35%  "
$ws.Range("D10").Value = "0.103"
$ws.Range("E10").Value = "  +5.83%  "
$ws.Range("D11").Value = "32.18"
$ws.Range("E11").Value = "  +21.28%  "
$ws.Range("D12").Value = "7.51"
$ws.Range("E12").Value = "  +20.18%  "
$ws.Range("E13").Value = "  +2.29%  "
$ws.Range("D14").Value = "2.712.95"
$ws.Range("E14").Value = "  -0.43%  "
$ws.Range("D15").Value = "16.93"
$ws.Range("E15").Value = "  +7.26%  "
$ws.Range("D16").Value = "0.916"
$ws.Range("E16").Value = "  +6.80%  "
$ws.Range("D17").Value = "2.359.56"
$ws.Range("E17").Value = "  -0.75%  "
$ws.Range("D18").Value = "44.450.69"
$ws.Range("E18").Value = "  +2.36%  "
$ws.Range("E19").Value = "  +4.66%  "
$ws.Range("D20").Value = "6.82"
$ws.Range("E20").Value = "  +7.36%  "
$ws.Range("D21").Value = "78.45"
$ws.Range("E21").Value = "  +5.67%  "
$ws.Range("D22").Value = "256.48"
$ws.Range("E22").Value = "  +1.94%  "
$ws.Range("B23").Value = "PancakeSwap"
$ws.Range("C23").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D23").Value = "2.59"
$ws.Range("E23").Value = "  +4.77%  "
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("E25").Value = "  -5.28%  "
$ws.Range("D26").Value = "10.78"
$ws.Range("E26").Value = "  +7.26%  "
$ws.Range("D27").Value = "2.31"
$ws.Range("E27").Value = "  +3.82%  "
$ws.Range("D28").Value = "22.62"
$ws.Range("E28").Value = "  -1.85%  "
$ws.Range("D29").Value = "175.12"
$ws.Range("E29").Value = "  +0.45%  "
$ws.Range("D30").Value = "1.60"
$ws.Range("E30").Value = "  +3.49%  "
$ws.Range("D31").Value = "0.131"
$ws.Range("E31").Value = "  +3.65%  "
$ws.Range("E32").Value = "  +4.66%  "
$ws.Range("D33").Value = "5.43"
$ws.Range("E33").Value = "  +8.55%  "
$ws.Range("E34").Value = "  +9.23%  "
$ws.Range("E35").Value = "  +5.15%  "
$ws.Range("D36").Value = "3.89"
$ws.Range("E36").Value = "  +6.75%  "
$ws.Range("D37").Value = "2.47"
$ws.Range("E37").Value = "  +0.55%  "
$ws.Range("D38").Value = "6.59"
$ws.Range("E38").Value = "  -0.29%  "
$ws.Range("E39").Value = "  +8.09%  "
$ws.Range("D40").Value = "19.30"
$ws.Range("E40").Value = "  +3.86%  "
$ws.Range("D41").Value = "8.99"
$ws.Range("E41").Value = "  +0.18%  "
$ws.Range("E42").Value = "  +0.17%  "
$ws.Range("E43").Value = "  +15.03%  "
$ws.Range("E44").Value = "  +3.15%  "
$ws.Range("E45").Value = "  +11.23%  "
$ws.Range("D46").Value = "0.0999"
$ws.Range("E46").Value = "  +5.51%  "
$ws.Range("D47").Value = "101.21"
$ws.Range("E47").Value = "  +1.34%  "
$ws.Range("D48").Value = "1.18"
$ws.Range("E48").Value = "  -0.60%  "
$ws.Range("E49").Value = "  -0.78%  "
$ws.Range("D50").Value = "1.455.24"
$ws.Range("E50").Value = "  +0.13%  "
$ws.Range("D51").Value = "0.000208"
$ws.Range("E51").Value = "  +3.18%  "
